$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 ("Hot Pursuit" overview, sldId 256) ------------------------------
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# 1) Reposition the existing "CNEOS Scout" flowchart box (shape id 10) -
#    it moves up and slightly to the right to make room for a new box below.
$scout = $s1.Shapes.Item("Flowchart: Alternate Process 9")
$scout.Left = 42.7847
$scout.Top = 199.6118
$scout.Width = 137.1177
$scout.Height = 35.4706

# 2) Add a new "NASA Horizons" flowchart box just below the relocated one.
#    Duplicating the existing "CNEOS Scout" shape keeps the same accent6
#    flowchart style (fill/line/effect/font refs) instead of the default.
$nasaRange = $scout.Duplicate()
$nasa = $nasaRange.Item(1)
$nasa.Name = "Flowchart: Alternate Process 20"
$nasa.Left = 43.6088
$nasa.Top = 279.3981
$nasa.Width = 137.1177
$nasa.Height = 35.4706
$nasa.TextFrame.TextRange.Text = "NASA Horizons"

# ---------------------------------------------------------------------------
# Slide 3 ("Site Characterization" process, sldId 259) ---------------------
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# 3) The "CNEOS / Scout" terminator shape (shape id 8) gains an "or / NASA /
#    Horizons" alternative underneath the existing two lines.
$cneos = $s3.Shapes.Item("Flowchart: Terminator 7")
$cneos.TextFrame.TextRange.Text = "CNEOS`rScout`ror`rNASA`rHorizons"

# 4) The "Site Translation" alternate-process box (shape id 22) gets a
#    second line clarifying it's conditional.
$siteTrans = $s3.Shapes.Item("Flowchart: Alternate Process 21")
$siteTrans.TextFrame.TextRange.Text = "Site Translation`r(if required)"
